$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2110
$ws.Range("I40").Value = 2000
$ws.Range("J40").Value = 2366.6667
$ws.Range("K40").Value = 2000
$ws.Range("L40").Value = 2366.6667
$ws.Range("M40").Value = -1825
$ws.Range("N40").Value = -2716.6667

$ws.Range("H64").Value = 3627.375
$ws.Range("I64").Value = 3698
$ws.Range("J64").Value = 3509.6667
$ws.Range("K64").Value = 3698
$ws.Range("L64").Value = 3509.6667
$ws.Range("M64").Value = -3450
$ws.Range("N64").Value = -4005.6667

$ws.Range("H67").Value = 3627.375
$ws.Range("I67").Value = 3698
$ws.Range("J67").Value = 3509.6667
$ws.Range("K67").Value = 3698
$ws.Range("L67").Value = 3509.6667
$ws.Range("M67").Value = -2840
$ws.Range("N67").Value = -5225.6667

$ws.Range("H80").Value = 1069.6923
$ws.Range("I80").Value = 1436.4
$ws.Range("J80").Value = 840.5
$ws.Range("K80").Value = 4309.200000000001
$ws.Range("L80").Value = 2521.5
$ws.Range("M80").Value = -3311.200000000001
$ws.Range("N80").Value = -4517.5

$ws.Range("H83").Value = 1069.6923
$ws.Range("I83").Value = 1436.4
$ws.Range("J83").Value = 840.5
$ws.Range("K83").Value = 12927.6
$ws.Range("L83").Value = 7564.5
$ws.Range("M83").Value = -7935.6
$ws.Range("N83").Value = -17548.5

$ws.Range("H86").Value = 4375
$ws.Range("I86").Value = 4457.143
$ws.Range("J86").Value = 3800
$ws.Range("K86").Value = 4457.143
$ws.Range("L86").Value = 3800
$ws.Range("M86").Value = -3334.143
$ws.Range("N86").Value = -6046

$ws.Range("H87").Value = 41504.2
$ws.Range("J87").Value = 41504.2
$ws.Range("L87").Value = 41504.2
$ws.Range("N87").Value = -44000.2

$ws.Range("H89").Value = 4375
$ws.Range("I89").Value = 4457.143
$ws.Range("J89").Value = 3800
$ws.Range("K89").Value = 22285.715
$ws.Range("L89").Value = 19000
$ws.Range("M89").Value = -16669.715
$ws.Range("N89").Value = -30232

$ws.Range("H90").Value = 41504.2
$ws.Range("J90").Value = 41504.2
$ws.Range("L90").Value = 124512.6
$ws.Range("N90").Value = -136992.6

$ws.Range("H116").Value = 5335.3335
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 5335.3335
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 5335.3335
$ws.Range("M116").Value = ""
$ws.Range("N116").Value = -12219.3335

$ws.Range("H129").Value = 710.17645
$ws.Range("I129").Value = 516.125
$ws.Range("J129").Value = 882.6667
$ws.Range("K129").Value = 1548.375
$ws.Range("L129").Value = 2648.0001
$ws.Range("M129").Value = 3451.625
$ws.Range("N129").Value = -12648.0001

$ws.Range("H138").Value = 1593.0862
$ws.Range("I138").Value = 1071.4615
$ws.Range("J138").Value = 1743.7778
$ws.Range("K138").Value = 3214.3845
$ws.Range("L138").Value = 5231.3334
$ws.Range("M138").Value = 1925.6155
$ws.Range("N138").Value = -15511.3334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H51").Value = 12000
$ws.Range("J51").Value = 12000
$ws.Range("L51").Value = 12000
$ws.Range("N51").Value = -13512

$ws.Range("H54").Value = 11000
$ws.Range("J54").Value = 11000
$ws.Range("L54").Value = 11000
$ws.Range("N54").Value = -12538

$ws.Range("H63").Value = 2062.7112
$ws.Range("I63").Value = 1929.0667
$ws.Range("K63").Value = 1929.0667
$ws.Range("M63").Value = -1243.0667

$ws.Range("H66").Value = 2062.7112
$ws.Range("I66").Value = 1929.0667
$ws.Range("K66").Value = 9645.333500000001
$ws.Range("M66").Value = -6213.333500000001

$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").Value = ""

$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").Value = ""

$ws.Range("H92").Value = 2503500
$ws.Range("J92").Value = 2503500
$ws.Range("L92").Value = 2503500
$ws.Range("N92").Value = -2508492

$ws.Range("H132").Value = 2790.9033
$ws.Range("I132").Value = 2456.2173
$ws.Range("J132").Value = 3753.125
$ws.Range("K132").Value = 7368.651899999999
$ws.Range("L132").Value = 11259.375
$ws.Range("M132").Value = -4838.651899999999
$ws.Range("N132").Value = -16319.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 491.57895
$ws.Range("J80").Value = 566.5625
$ws.Range("L80").Value = 566.5625
$ws.Range("N80").Value = -2562.5625

$ws.Range("H83").Value = 491.57895
$ws.Range("J83").Value = 566.5625
$ws.Range("L83").Value = 2832.8125
$ws.Range("N83").Value = -12816.8125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 7409637
$ws.Range("I62").Value = 2347.6191
$ws.Range("J62").Value = 33335150
$ws.Range("K62").Value = 2347.6191
$ws.Range("L62").Value = 33335150
$ws.Range("M62").Value = -1723.6191
$ws.Range("N62").Value = -33336398

$ws.Range("H65").Value = 7409637
$ws.Range("I65").Value = 2347.6191
$ws.Range("J65").Value = 33335150
$ws.Range("K65").Value = 11738.0955
$ws.Range("L65").Value = 166675750
$ws.Range("M65").Value = -8618.095499999999
$ws.Range("N65").Value = -166681990

$ws.Range("H132").Value = 2092.8
$ws.Range("I132").Value = 1821.2222
$ws.Range("K132").Value = 5463.6666
$ws.Range("M132").Value = -2933.6666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 350
$ws.Range("I15").Value = 350
$ws.Range("K15").Value = 1050
$ws.Range("M15").Value = -910

$ws.Range("H20").Value = 287.5
$ws.Range("J20").Value = 275
$ws.Range("L20").Value = 825
$ws.Range("N20").Value = -1279

$ws.Range("H26").Value = 275
$ws.Range("J26").Value = 293.33334
$ws.Range("L26").Value = 880.0000200000001
$ws.Range("N26").Value = -1456.00002

$ws.Range("H32").Value = 1920.2
$ws.Range("J32").Value = 2262.5
$ws.Range("L32").Value = 6787.5
$ws.Range("N32").Value = -7353.5

$ws.Range("H131").Value = 22226168
$ws.Range("J131").Value = 4724.757
$ws.Range("L131").Value = 14174.271
$ws.Range("N131").Value = -24254.271

$ws.Range("H139").Value = 1807.5151
$ws.Range("I139").Value = 1942.4
$ws.Range("K139").Value = 5827.200000000001
$ws.Range("M139").Value = -687.2000000000007

$ws.Range("H140").Value = 28123.525
$ws.Range("J140").Value = 3025.2415
$ws.Range("L140").Value = 9075.7245
$ws.Range("N140").Value = -19435.7245

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 811.5833
$ws.Range("I102").Value = 762.5
$ws.Range("K102").Value = 762.5
$ws.Range("M102").Value = 859.5

$ws.Range("H122").Value = 2045
$ws.Range("I122").Value = 2110.4167
$ws.Range("J122").Value = 1932.8572
$ws.Range("K122").Value = 6331.250100000001
$ws.Range("L122").Value = 5798.571599999999
$ws.Range("M122").Value = -3881.250100000001
$ws.Range("N122").Value = -10698.5716

$ws.Range("H132").Value = 3759.2222
$ws.Range("I132").Value = 3618.1428
$ws.Range("J132").Value = 4253
$ws.Range("K132").Value = 10854.4284
$ws.Range("L132").Value = 12759
$ws.Range("M132").Value = -8324.428400000001
$ws.Range("N132").Value = -17819

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 62501788
$ws.Range("I122").Value = 125000750
$ws.Range("K122").Value = 375002250
$ws.Range("M122").Value = -374999800

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 563.8125
$ws.Range("I107").Value = 557.2857
$ws.Range("J107").Value = 568.8889
$ws.Range("K107").Value = 1671.8571
$ws.Range("L107").Value = 1706.6667
$ws.Range("M107").Value = 248.1428999999998
$ws.Range("N107").Value = -5546.6667

$ws.Range("H136").Value = 995.0323
$ws.Range("I136").Value = 930.1429000000001
$ws.Range("J136").Value = 1131.3
$ws.Range("K136").Value = 2790.4287
$ws.Range("L136").Value = 3393.9
$ws.Range("M136").Value = -240.4287000000004
$ws.Range("N136").Value = -8493.9
